# NATMI LR-pairs output refresh (Efnb3-Ephb3, YoungD7) with new TPM-based results.
# The underlying analysis now reports all three clusters (ECs, FAPs, MuSCs) as both
# sending and target clusters (a full 3x3 = 9-row grid in rows 2-10), and all of the
# expression / specificity statistics were recomputed against the new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (ligand Efnb3 / receptor Ephb3)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb3"
$ws.Range("C2").Value = "Ephb3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2151756666666667
$ws.Range("H2").Value = 0.645527
$ws.Range("I2").Value = 0.1791915537270423
$ws.Range("J2").Value = 0.1791915537270423
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1709536666666667
$ws.Range("N2").Value = 0.512861
$ws.Range("O2").Value = 0.007882947722998253
$ws.Range("P2").Value = 0.007882947722998253
$ws.Range("Q2").Value = 0.03678506919411111
$ws.Range("R2").Value = 0.331065622747
$ws.Range("S2").Value = 0.001412557650433107
$ws.Range("T2").Value = 0.001412557650433107

# Row 3: ECs -> FAPs (ligand Efnb3 / receptor Ephb3)
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efnb3"
$ws.Range("C3").Value = "Ephb3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2151756666666667
$ws.Range("H3").Value = 0.645527
$ws.Range("I3").Value = 0.1791915537270423
$ws.Range("J3").Value = 0.1791915537270423
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 17.192962
$ws.Range("N3").Value = 51.578886
$ws.Range("O3").Value = 0.7927950496303802
$ws.Range("P3").Value = 0.7927950496303802
$ws.Range("Q3").Value = 3.699507060324666
$ws.Range("R3").Value = 33.295563542922
$ws.Range("S3").Value = 0.1420621767303754
$ws.Range("T3").Value = 0.1420621767303754

# Row 4: ECs -> MuSCs (ligand Efnb3 / receptor Ephb3)
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efnb3"
$ws.Range("C4").Value = "Ephb3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2151756666666667
$ws.Range("H4").Value = 0.645527
$ws.Range("I4").Value = 0.1791915537270423
$ws.Range("J4").Value = 0.1791915537270423
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.322599666666666
$ws.Range("N4").Value = 12.967799
$ws.Range("O4").Value = 0.1993220026466216
$ws.Range("P4").Value = 0.1993220026466216
$ws.Range("Q4").Value = 0.9301182650081109
$ws.Range("R4").Value = 8.371064385073
$ws.Range("S4").Value = 0.03571681934623376
$ws.Range("T4").Value = 0.03571681934623376

# Row 5: FAPs -> ECs (ligand Efnb3 / receptor Ephb3)
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efnb3"
$ws.Range("C5").Value = "Ephb3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.5010686666666667
$ws.Range("H5").Value = 1.503206
$ws.Range("I5").Value = 0.4172742870736815
$ws.Range("J5").Value = 0.4172742870736815
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1709536666666667
$ws.Range("N5").Value = 0.512861
$ws.Range("O5").Value = 0.007882947722998253
$ws.Range("P5").Value = 0.007882947722998253
$ws.Range("Q5").Value = 0.08565952581844445
$ws.Range("R5").Value = 0.770935732366
$ws.Range("S5").Value = 0.003289351391153196
$ws.Range("T5").Value = 0.003289351391153196

# Row 6: FAPs -> FAPs (ligand Efnb3 / receptor Ephb3)
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efnb3"
$ws.Range("C6").Value = "Ephb3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.5010686666666667
$ws.Range("H6").Value = 1.503206
$ws.Range("I6").Value = 0.4172742870736815
$ws.Range("J6").Value = 0.4172742870736815
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 17.192962
$ws.Range("N6").Value = 51.578886
$ws.Range("O6").Value = 0.7927950496303802
$ws.Range("P6").Value = 0.7927950496303802
$ws.Range("Q6").Value = 8.614854545390667
$ws.Range("R6").Value = 77.533690908516
$ws.Range("S6").Value = 0.3308129891300608
$ws.Range("T6").Value = 0.3308129891300608

# Row 7: FAPs -> MuSCs (ligand Efnb3 / receptor Ephb3)
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efnb3"
$ws.Range("C7").Value = "Ephb3"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.5010686666666667
$ws.Range("H7").Value = 1.503206
$ws.Range("I7").Value = 0.4172742870736815
$ws.Range("J7").Value = 0.4172742870736815
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.322599666666666
$ws.Range("N7").Value = 12.967799
$ws.Range("O7").Value = 0.1993220026466216
$ws.Range("P7").Value = 0.1993220026466216
$ws.Range("Q7").Value = 2.165919251510445
$ws.Range("R7").Value = 19.493273263594
$ws.Range("S7").Value = 0.08317194655246749
$ws.Range("T7").Value = 0.08317194655246749

# Row 8: MuSCs -> ECs (ligand Efnb3 / receptor Ephb3)
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Efnb3"
$ws.Range("C8").Value = "Ephb3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.4845693333333334
$ws.Range("H8").Value = 1.453708
$ws.Range("I8").Value = 0.4035341591992763
$ws.Range("J8").Value = 0.4035341591992763
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1709536666666667
$ws.Range("N8").Value = 0.512861
$ws.Range("O8").Value = 0.007882947722998253
$ws.Range("P8").Value = 0.007882947722998253
$ws.Range("Q8").Value = 0.08283890428755555
$ws.Range("R8").Value = 0.745550138588
$ws.Range("S8").Value = 0.003181038681411949
$ws.Range("T8").Value = 0.003181038681411949

# Row 9: MuSCs -> FAPs (ligand Efnb3 / receptor Ephb3)
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Efnb3"
$ws.Range("C9").Value = "Ephb3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.4845693333333334
$ws.Range("H9").Value = 1.453708
$ws.Range("I9").Value = 0.4035341591992763
$ws.Range("J9").Value = 0.4035341591992763
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 17.192962
$ws.Range("N9").Value = 51.578886
$ws.Range("O9").Value = 0.7927950496303802
$ws.Range("P9").Value = 0.7927950496303802
$ws.Range("Q9").Value = 8.331182134365333
$ws.Range("R9").Value = 74.980639209288
$ws.Range("S9").Value = 0.319919883769944
$ws.Range("T9").Value = 0.319919883769944

# Row 10: MuSCs -> MuSCs (ligand Efnb3 / receptor Ephb3)
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Efnb3"
$ws.Range("C10").Value = "Ephb3"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.4845693333333334
$ws.Range("H10").Value = 1.453708
$ws.Range("I10").Value = 0.4035341591992763
$ws.Range("J10").Value = 0.4035341591992763
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 4.322599666666666
$ws.Range("N10").Value = 12.967799
$ws.Range("O10").Value = 0.1993220026466216
$ws.Range("P10").Value = 0.1993220026466216
$ws.Range("Q10").Value = 2.094599238743555
$ws.Range("R10").Value = 18.851393148692
$ws.Range("S10").Value = 0.08043323674792037
$ws.Range("T10").Value = 0.08043323674792037
